$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: fix the "Day 8, 9" date text (4-5 -> 04-05); B8 keeps the
#     "background-blend-mode" text (unchanged content, just shifts which
#     shared-string slot it occupies once a new row/string is introduced).
$ws.Range("A8").Value2 = "Day 8, 9 (04-05/07/2019"
$ws.Range("B8").Value2 = "Learn about new CSS brand feature: background-blend-mode, box-decoration-break"

# --- Row 9: fix the "Day 10" date text (6 -> 06); B9 keeps the
#     "over-flow" text.
$ws.Range("A9").Value2 = "Day 10 (06/07/2019"
$ws.Range("B9").Value2 = "Learn how to use over-flow: hidden when before we used clip-path"

# --- Row 10 (new row): "Day 11" entry plus a rich-text note about the
#     shape-outside CSS property, with "shape-outside: " rendered bold.
$ws.Range("A10").Value2 = "Day 11 (08/07/2019"

# Copy column A's center/center alignment style down onto the new row
# (xlPasteFormats = -4122) so A10 matches A1:A9's look.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$cellB10 = $ws.Range("B10")
$noteText = "Learing new some properties: shape-outside: define where the content floats around the element => if you actually want the element to look lije that circle, we can then use the clip-path property"
$cellB10.Value2 = $noteText

# Make sure the bold font used for the run gets registered in the
# styles font table (mirrors how Excel records a distinct font entry
# for the bold "shape-outside: " run).
$cellB10.Font.Bold = $true
$cellB10.Font.Bold = $false

$boldPhrase = "shape-outside: "
$boldStart = $noteText.IndexOf($boldPhrase) + 1
$boldLen = $boldPhrase.Length
$leadLen = $boldStart - 1
$tailStart = $boldStart + $boldLen
$tailLen = $noteText.Length - ($leadLen + $boldLen)

# Three runs total: plain lead-in, bold "shape-outside: ", plain remainder.
# Re-asserting the (already-default) font on the plain runs forces this
# engine to emit an explicit run-level <rPr> for them too, matching the
# three-run rich-text layout Excel produced for this cell.
$cellB10.Characters(1, $leadLen).Font.Name = "Calibri"
$cellB10.Characters($boldStart, $boldLen).Font.Bold = $true
$cellB10.Characters($tailStart, $tailLen).Font.Name = "Calibri"

# --- Update the active selection left behind by the edit (A14, not a
#     real data cell -- just reflects where the editor's cursor ended up).
$ws.Range("A14").Select()
